$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames
$ws.Range("C1").Value = "rules"
$ws.Range("E1").Value = "adaptive_filter"

# Row data updates: E column becomes text "RLS", F/G/H updated with refined precision values
$ws.Range("E2").Value = "RLS"
$ws.Range("F2").Value = 903.6050430168957
$ws.Range("G2").Value = 1.969301855064028
$ws.Range("H2").Value = 715.1660566772258

$ws.Range("E3").Value = "RLS"
$ws.Range("F3").Value = 784.4210104514062
$ws.Range("G3").Value = 1.709554149759509
$ws.Range("H3").Value = 626.0312675942313

$ws.Range("E4").Value = "RLS"
$ws.Range("F4").Value = 484.6198277468936
$ws.Range("G4").Value = 1.056172420858128
$ws.Range("H4").Value = 396.8366198010475

$ws.Range("E5").Value = "RLS"
$ws.Range("F5").Value = 298.7602704863381
$ws.Range("G5").Value = 0.6511131820644072
$ws.Range("H5").Value = 241.0604579601335

$ws.Range("E6").Value = "RLS"
$ws.Range("F6").Value = 181.5948014225014
$ws.Range("G6").Value = 0.3957647006011996
$ws.Range("H6").Value = 145.5795370904831

$ws.Range("E7").Value = "RLS"
$ws.Range("F7").Value = 176.4698201245835
$ws.Range("G7").Value = 0.3845954012981971
$ws.Range("H7").Value = 140.143364525665

$ws.Range("E8").Value = "RLS"
$ws.Range("F8").Value = 173.1408329574672
$ws.Range("G8").Value = 0.3773402618383747
$ws.Range("H8").Value = 136.1794149930056

$ws.Range("E9").Value = "RLS"
$ws.Range("F9").Value = 171.0019250903222
$ws.Range("G9").Value = 0.3726787614814085
$ws.Range("H9").Value = 132.9647137220902

$ws.Range("E10").Value = "RLS"
$ws.Range("F10").Value = 169.3699062348481
$ws.Range("G10").Value = 0.3691219666356714
$ws.Range("H10").Value = 130.7663926344838

$ws.Range("E11").Value = "RLS"
$ws.Range("F11").Value = 167.5922069020842
$ws.Range("G11").Value = 0.3652476781721301
$ws.Range("H11").Value = 129.272482429537

$ws.Range("E12").Value = "RLS"
$ws.Range("F12").Value = 165.2361489177383
$ws.Range("G12").Value = 0.3601129244486231
$ws.Range("H12").Value = 127.3216455569816

$ws.Range("E13").Value = "RLS"
$ws.Range("F13").Value = 162.4936663865674
$ws.Range("G13").Value = 0.3541360034720825
$ws.Range("H13").Value = 125.149027249837

$ws.Range("E14").Value = "RLS"
$ws.Range("F14").Value = 160.5467869559622
$ws.Range("G14").Value = 0.3498930067071727
$ws.Range("H14").Value = 124.0260370738752

$ws.Range("E15").Value = "RLS"
$ws.Range("F15").Value = 160.3422495722593
$ws.Range("G15").Value = 0.349447241322983
$ws.Range("H15").Value = 125.1036273762254

$ws.Range("E16").Value = "RLS"
$ws.Range("F16").Value = 163.0355370277461
$ws.Range("G16").Value = 0.3553169473669009
$ws.Range("H16").Value = 127.9919714595157
